$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BL1:BL11").Copy()
$ws.Range("BM1:BM11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("BM1").Value = "08-sep"

$ws.Range("P2").Copy()
$ws.Range("BM2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("BM2").Value = 16
$ws.Range("BM3").Value = 16
$ws.Range("BM4").Value = 12
$ws.Range("BM5").Value = 12
$ws.Range("BM6").Value = 9
$ws.Range("BM7").Value = 17
$ws.Range("BM8").Value = 23
$ws.Range("BM9").Value = 14
$ws.Range("BM10").Value = 13
$ws.Range("BM11").Value = 12

$ws.Range("BM2").Select()
